# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures on the
# crypto symbol list, as produced by the Thu Jan 26 05:37:44 UTC 2023
# GitHub Actions scrape. Both columns hold plain text in the workbook
# (e.g. "307.33", "1.64%"), so each value is written with a leading
# apostrophe to force Excel to keep it as literal text instead of
# auto-converting it to a number/percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.33"
$ws.Range("E2").Value = "'1.64%"
$ws.Range("D3").Value = "'36.20"
$ws.Range("E3").Value = "'3.08%"
$ws.Range("D4").Value = "'5.101"
$ws.Range("E4").Value = "'1.05%"
$ws.Range("D5").Value = "'0.08126"
$ws.Range("E5").Value = "'2.79%"
$ws.Range("D6").Value = "'1.940"
$ws.Range("E6").Value = "'-0.09%"
$ws.Range("D7").Value = "'4.186"
$ws.Range("E7").Value = "'4.12%"
$ws.Range("D8").Value = "'7.787"
$ws.Range("E8").Value = "'0.86%"
$ws.Range("D9").Value = "'0.9308"
$ws.Range("E9").Value = "'0.76%"
$ws.Range("D10").Value = "'0.1400"
$ws.Range("E10").Value = "'17.65%"
$ws.Range("D11").Value = "'0.1921"
$ws.Range("E11").Value = "'4.52%"
$ws.Range("D12").Value = "'0.09213"
$ws.Range("E12").Value = "'-1.81%"
$ws.Range("D13").Value = "'0.03418"
$ws.Range("E13").Value = "'-3.35%"
$ws.Range("D14").Value = "'0.09862"
$ws.Range("E14").Value = "'-0.28%"
$ws.Range("D15").Value = "'0.001417"
$ws.Range("E15").Value = "'2.13%"
$ws.Range("D16").Value = "'0.005735"
$ws.Range("E16").Value = "'-1.47%"
$ws.Range("D17").Value = "'3.606"
$ws.Range("E17").Value = "'3.28%"
$ws.Range("D18").Value = "'2.972"
$ws.Range("D19").Value = "'0.3439"
$ws.Range("E19").Value = "'-0.14%"
$ws.Range("D21").Value = "'4.894"
$ws.Range("E21").Value = "'-2.83%"
$ws.Range("D23").Value = "'0.04512"
$ws.Range("E23").Value = "'0.56%"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'0.09%"
$ws.Range("D25").Value = "'0.004864"
$ws.Range("E25").Value = "'6.41%"
$ws.Range("D26").Value = "'0.0001241"
$ws.Range("E26").Value = "'-0.70%"
$ws.Range("D39").Value = "'0.02002"
$ws.Range("E39").Value = "'5.08%"
$ws.Range("D40").Value = "'0.04943"
$ws.Range("E40").Value = "'5.10%"
$ws.Range("D41").Value = "'0.007655"
$ws.Range("E41").Value = "'0.77%"
$ws.Range("D42").Value = "'0.01023"
$ws.Range("E42").Value = "'7.07%"
$ws.Range("D43").Value = "'0.1381"
$ws.Range("E43").Value = "'4.38%"
$ws.Range("D44").Value = "'0.002102"
$ws.Range("E44").Value = "'-0.38%"
$ws.Range("D45").Value = "'0.01154"
$ws.Range("E45").Value = "'2.97%"
$ws.Range("D46").Value = "'0.00006459"
$ws.Range("E46").Value = "'7.29%"
$ws.Range("E47").Value = "'0.05%"
$ws.Range("E49").Value = "'-8.67%"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E51").Value = "'0.05%"
